$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's used range started at row 2 (A2:I26) - row 1 was blank/unused.
# Remove that blank leading row so every record shifts up by one
# (old row 2 -> row 1, old row 3 -> row 2, ... old row 26 -> row 25).
$ws.Rows("1").Delete()

# Label the now-empty header cell that used to hold the year columns.
$ws.Range("A1").Value = "Mes"

# Leave the selection where the author left it after the edit.
[void]$ws.Range("C23").Select()
